# Adds a "Children" column (L) to the Employees export sheet, mirroring
# the "Spouse" column that already exists, and leaves the selection on
# the newly-widened "Nationality" column (E) as in the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell L1 -> shared string "Children" (goes in right after
# the existing "Spouse" header in K1).
$ws.Range("L1").Value = "Children"

# The source workbook also picked up an explicit width for column E
# ("Nationality") -- set it as close as the host can represent.
$ws.Columns.Item(5).ColumnWidth = 10

# Selection moves to E1 in the edited workbook.
$ws.Range("E1").Select() | Out-Null
